# Auto-generated edit script: updates crypto price/volume table to reflect
# the latest GitHub Actions scrape, incl. the Aave/WhiteBITCoin row swap.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.939.06"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "2.500.13"
$ws.Range("E3").Value = "  -0.84%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'535.52"
$ws.Range("E5").Value = "  -0.32%  "
$ws.Range("D6").Value = "'138.00"
$ws.Range("E6").Value = "  -1.56%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +0.79%  "
$ws.Range("D9").Value = "2.524.37"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  +0.93%  "
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'5.37"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "2.974.06"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").Value = "'23.26"
$ws.Range("E15").Value = "  +1.50%  "
$ws.Range("D16").Value = "58.919.18"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "2.512.43"
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").Value = "'11.08"
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").Value = "'4.26"
$ws.Range("E20").Value = "  +0.94%  "
$ws.Range("D21").Value = "'325.51"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").Value = "'64.75"
$ws.Range("E24").Value = "  +4.23%  "
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("D29").Value = "0.0₃0777"
$ws.Range("E29").Value = "  +1.30%  "
$ws.Range("D30").Value = "'6.72"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("D31").Value = "'1.77"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").Value = "'167.91"
$ws.Range("E32").Value = "  +4.15%  "
$ws.Range("D33").Value = "'1.17"
$ws.Range("E33").Value = "  +4.35%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "'1.40"
$ws.Range("E35").Value = "  -4.37%  "
$ws.Range("D36").Value = "'18.57"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "'4.13"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("D38").Value = "'1.56"
$ws.Range("E38").Value = "  -1.10%  "
$ws.Range("D39").Value = "'36.75"
$ws.Range("E39").Value = "  -0.60%  "
$ws.Range("D40").Value = "'0.833"
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("E41").Value = "  -0.17%  "
$ws.Range("D42").Value = "'5.31"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("D43").Value = "'282.72"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").Value = "'0.994"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "'0.605"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'129.28"
$ws.Range("E46").Value = "  +5.62%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'10.88"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("D48").Value = "'0.0932"
$ws.Range("E48").Value = "  +0.15%  "
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  +0.46%  "
$ws.Range("D51").Value = "'17.39"
$ws.Range("E51").Value = "  -0.29%  "
